# Fill in a new guest row (row 8) in the GuestList table and move the
# active cell selection, matching the author's manual edit in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "takale letikun"
$ws.Range("B8").Value = 6
$ws.Range("C8").Value = "groom"
$ws.Range("D8").Value = "school friends"

$ws.Range("B20").Select()
